# Updated legacy GSC export data
# - "Chart" sheet: rolling date window advanced by one day (oldest date
#   dropped, newest days appended) and all Not-indexed/Indexed/Impressions
#   figures re-pulled for the new window.
# - "Critical issues" sheet: "Crawled - currently not indexed" page count
#   bumped 3 -> 4.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Chart sheet: Date / Not indexed / Indexed / Impressions
# ---------------------------------------------------------------------
$chart = $wb.Worksheets.Item("Chart")

# Column A holds dates written as literal text (e.g. "2025-10-17"), not
# real Excel date serials, so format the column as Text first to stop
# Excel's autodetection from converting the strings into dates.
$chart.Range("A2:A90").NumberFormat = "@"

$rows = @(
    @(2, "2025-10-17", $null, $null, 60),
    @(3, "2025-10-18", 8, 127, 61),
    @(4, "2025-10-19", 9, 126, 54),
    @(5, "2025-10-20", 9, 126, 51),
    @(6, "2025-10-21", 9, 126, 36),
    @(7, "2025-10-22", 14, 121, 55),
    @(8, "2025-10-23", 14, 121, 41),
    @(9, "2025-10-24", 14, 121, 56),
    @(10, "2025-10-25", 14, 121, 15),
    @(11, "2025-10-26", 20, 115, 21),
    @(12, "2025-10-27", 20, 115, 24),
    @(13, "2025-10-28", 20, 115, 23),
    @(14, "2025-10-29", 21, 114, 28),
    @(15, "2025-10-30", 21, 114, 25),
    @(16, "2025-10-31", 21, 114, 30),
    @(17, "2025-11-01", 21, 114, 23),
    @(18, "2025-11-02", 28, 107, 25),
    @(19, "2025-11-03", 28, 107, 27),
    @(20, "2025-11-04", 28, 107, 28),
    @(21, "2025-11-05", 36, 99, 13),
    @(22, "2025-11-06", 36, 99, 19),
    @(23, "2025-11-07", 36, 99, 9),
    @(24, "2025-11-08", 36, 99, 10),
    @(25, "2025-11-09", 36, 99, 15),
    @(26, "2025-11-10", 36, 99, 13),
    @(27, "2025-11-11", 36, 99, 9),
    @(28, "2025-11-12", 38, 97, 16),
    @(29, "2025-11-13", 38, 97, 11),
    @(30, "2025-11-14", 38, 97, 10),
    @(31, "2025-11-15", 38, 97, 8),
    @(32, "2025-11-16", 43, 92, 18),
    @(33, "2025-11-17", 43, 92, 16),
    @(34, "2025-11-18", 43, 92, 20),
    @(35, "2025-11-19", 46, 89, 11),
    @(36, "2025-11-20", 46, 89, 15),
    @(37, "2025-11-21", 46, 89, 14),
    @(38, "2025-11-22", 46, 89, 14),
    @(39, "2025-11-23", 73, 62, 9),
    @(40, "2025-11-24", 73, 62, 5),
    @(41, "2025-11-25", 73, 62, 14),
    @(42, "2025-11-26", 73, 62, 6),
    @(43, "2025-11-27", 73, 62, 2),
    @(44, "2025-11-28", 73, 62, 2),
    @(45, "2025-11-29", 73, 62, 0),
    @(46, "2025-11-30", 73, 62, 7),
    @(47, "2025-12-01", 73, 62, 5),
    @(48, "2025-12-02", 73, 62, 4),
    @(49, "2025-12-03", 73, 62, 1),
    @(50, "2025-12-04", 73, 62, 6),
    @(51, "2025-12-05", 73, 62, 2),
    @(52, "2025-12-06", 73, 62, 2),
    @(53, "2025-12-07", 73, 62, 4),
    @(54, "2025-12-08", 73, 62, 3),
    @(55, "2025-12-09", 73, 62, 2),
    @(56, "2025-12-10", 73, 62, 2),
    @(57, "2025-12-11", 73, 62, 2),
    @(58, "2025-12-12", 73, 62, 3),
    @(59, "2025-12-13", 73, 62, 3),
    @(60, "2025-12-14", 73, 62, 16),
    @(61, "2025-12-15", 73, 62, 4),
    @(62, "2025-12-16", 79, 56, 9),
    @(63, "2025-12-17", 79, 56, 9),
    @(64, "2025-12-18", 79, 56, 15),
    @(65, "2025-12-19", 79, 56, 7),
    @(66, "2025-12-20", 79, 56, 7),
    @(67, "2025-12-21", 79, 56, 3),
    @(68, "2025-12-22", 79, 56, 12),
    @(69, "2025-12-23", 79, 56, 2),
    @(70, "2025-12-24", 85, 50, 1),
    @(71, "2025-12-25", 85, 50, 8),
    @(72, "2025-12-26", 85, 50, 2),
    @(73, "2025-12-27", 85, 50, 2),
    @(74, "2025-12-28", 85, 50, 2),
    @(75, "2025-12-29", 85, 50, 0),
    @(76, "2025-12-30", 85, 50, 1),
    @(77, "2025-12-31", 85, 50, 72),
    @(78, "2026-01-01", 85, 50, 1),
    @(79, "2026-01-02", 85, 50, 2),
    @(80, "2026-01-03", 85, 50, 29),
    @(81, "2026-01-04", 90, 45, 10),
    @(82, "2026-01-05", 90, 45, 2),
    @(83, "2026-01-06", 90, 45, 29),
    @(84, "2026-01-07", 93, 42, 40),
    @(85, "2026-01-08", 93, 42, 14),
    @(86, "2026-01-09", 93, 42, 85),
    @(87, "2026-01-10", 93, 42, 121),
    @(88, "2026-01-11", 94, 41, 9),
    @(89, "2026-01-12", 94, 41, 26),
    @(90, "2026-01-13", 94, 41, 67)
)

foreach ($row in $rows) {
    $r = $row[0]
    $chart.Cells.Item($r, 1).Value = $row[1]
    if ($null -ne $row[2]) {
        $chart.Cells.Item($r, 2).Value = $row[2]
    }
    if ($null -ne $row[3]) {
        $chart.Cells.Item($r, 3).Value = $row[3]
    }
    $chart.Cells.Item($r, 4).Value = $row[4]
}

# ---------------------------------------------------------------------
# 2) Critical issues sheet: updated "Crawled - currently not indexed"
#    page count (Source=Google systems, Validation=Passed): 3 -> 4
# ---------------------------------------------------------------------
$critical = $wb.Worksheets.Item("Critical issues")
$critical.Cells.Item(6, 4).Value = 4
